$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 20, shifting existing rows 20-27 down to 21-28
$ws.Rows.Item(20).Insert()

# Copy the style of the date cell from the row below (now row 21) into the new row 20, D column
$ws.Range("D21").Copy()
$ws.Range("D20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 20 with the new weekly entry
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = 45146
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 100112035
$ws.Range("G20").Value = "Bruselas (repollito)"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 140
$ws.Range("K20").Value = 26000
$ws.Range("L20").Value = 26000
$ws.Range("M20").Value = 26000
$ws.Range("N20").Value = "`$/malla 15 kilos"
$ws.Range("O20").Value = "Provincia de Quillota"
$ws.Range("P20").Value = 1733
$ws.Range("Q20").Value = 15
$ws.Range("R20").Value = "Hortaliza"
